# Apply the edit described by the diff:
#  - Insert 3 new rows before row 483 (shifting existing rows 483-545 down to 486-548)
#  - Populate the 3 new rows (483-485) with new "Valencia" price data
#
# Resulting sheet dimension becomes A1:T548 (was A1:T545).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 blank rows starting at row 483 ------------------------------
$ws.Rows("483:485").Insert()

# --- 2. Fill in the 3 new rows with their full data --------------------------
# Common (constant) metadata columns shared by every row in this block.
$marketId   = 7
$marketName = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$tipo       = "Fruta"
$productId  = 100102
$producto   = "Cítricos"
$catId      = 100102005
$categoria  = "Naranja"
$origen     = "Región de O'Higgins"
$unidad     = "$/bandeja 15 kilos granel"

# Row 483: Valencia / Especial
$ws.Cells.Item(483, 1).Value  = $marketId
$ws.Cells.Item(483, 2).Value  = $marketName
$ws.Cells.Item(483, 3).Value  = $region
$ws.Cells.Item(483, 4).Value  = 44984
$ws.Cells.Item(483, 5).Value  = $codreg
$ws.Cells.Item(483, 6).Value  = $tipo
$ws.Cells.Item(483, 7).Value  = $productId
$ws.Cells.Item(483, 8).Value  = $producto
$ws.Cells.Item(483, 9).Value  = $catId
$ws.Cells.Item(483, 10).Value = $categoria
$ws.Cells.Item(483, 11).Value = "Valencia"
$ws.Cells.Item(483, 12).Value = "Especial"
$ws.Cells.Item(483, 13).Value = 50
$ws.Cells.Item(483, 14).Value = 14000
$ws.Cells.Item(483, 15).Value = 14000
$ws.Cells.Item(483, 16).Value = 14000
$ws.Cells.Item(483, 17).Value = $unidad
$ws.Cells.Item(483, 18).Value = $origen
$ws.Cells.Item(483, 19).Value = 933
$ws.Cells.Item(483, 20).Value = 15

# Row 484: Valencia / Primera
$ws.Cells.Item(484, 1).Value  = $marketId
$ws.Cells.Item(484, 2).Value  = $marketName
$ws.Cells.Item(484, 3).Value  = $region
$ws.Cells.Item(484, 4).Value  = 44984
$ws.Cells.Item(484, 5).Value  = $codreg
$ws.Cells.Item(484, 6).Value  = $tipo
$ws.Cells.Item(484, 7).Value  = $productId
$ws.Cells.Item(484, 8).Value  = $producto
$ws.Cells.Item(484, 9).Value  = $catId
$ws.Cells.Item(484, 10).Value = $categoria
$ws.Cells.Item(484, 11).Value = "Valencia"
$ws.Cells.Item(484, 12).Value = "Primera"
$ws.Cells.Item(484, 13).Value = 50
$ws.Cells.Item(484, 14).Value = 12000
$ws.Cells.Item(484, 15).Value = 12000
$ws.Cells.Item(484, 16).Value = 12000
$ws.Cells.Item(484, 17).Value = $unidad
$ws.Cells.Item(484, 18).Value = $origen
$ws.Cells.Item(484, 19).Value = 800
$ws.Cells.Item(484, 20).Value = 15

# Row 485: Valencia / Segunda
$ws.Cells.Item(485, 1).Value  = $marketId
$ws.Cells.Item(485, 2).Value  = $marketName
$ws.Cells.Item(485, 3).Value  = $region
$ws.Cells.Item(485, 4).Value  = 44984
$ws.Cells.Item(485, 5).Value  = $codreg
$ws.Cells.Item(485, 6).Value  = $tipo
$ws.Cells.Item(485, 7).Value  = $productId
$ws.Cells.Item(485, 8).Value  = $producto
$ws.Cells.Item(485, 9).Value  = $catId
$ws.Cells.Item(485, 10).Value = $categoria
$ws.Cells.Item(485, 11).Value = "Valencia"
$ws.Cells.Item(485, 12).Value = "Segunda"
$ws.Cells.Item(485, 13).Value = 50
$ws.Cells.Item(485, 14).Value = 10000
$ws.Cells.Item(485, 15).Value = 10000
$ws.Cells.Item(485, 16).Value = 10000
$ws.Cells.Item(485, 17).Value = $unidad
$ws.Cells.Item(485, 18).Value = $origen
$ws.Cells.Item(485, 19).Value = 667
$ws.Cells.Item(485, 20).Value = 15

# --- 3. Make sure the date column keeps its date number format ---------------
$ws.Range("D483:D485").NumberFormat = $ws.Range("D486").NumberFormat
